# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (number of interested attendees) column (F) for a
# handful of rows on the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - first worksheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 6020
$wsExpo.Range("F6").Value = 113
$wsExpo.Range("F8").Value = 62
$wsExpo.Range("F9").Value = 554
$wsExpo.Range("F10").Value = 33

# Sheet "全部类型" (All Types) - fourth worksheet, same data duplicated with
# two extra rows shifted down by one (rows 5 and 9 inserted earlier).
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 6020
$wsAll.Range("F7").Value = 113
$wsAll.Range("F10").Value = 62
$wsAll.Range("F11").Value = 554
$wsAll.Range("F12").Value = 33
